# Updated cryptos list on Sat Oct 26 21:55:49 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table on
# Sheet1 with the latest scraped figures, and fixes the BabyDogeCoin/Cronos
# rows (50/51), which had swapped - row 50 should hold Cronos, row 51 should
# hold BabyDogeCoin (their rank numbers in column A stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many "Price" values (column D) are plain decimals ("7.39", "1.00", ...).
# Excel's normal Range.Value setter auto-detects those as numbers, which
# would turn "1.00" into 1 and drop the trailing zero. To keep them as the
# literal text the source sheet stores, force the cell to Text format only
# for the moment of the write, then clear that formatting again so the
# cell's style stays at its original (default) index.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Row 2: Bitcoin ---
Set-TextValue "D2" '67.105.50'
$ws.Range("E2").Value = '  -0.15%  '

# --- Row 3: Ethereum ---
Set-TextValue "D3" '2.485.87'

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = '  +0.03%  '

# --- Row 5: BNB ---
Set-TextValue "D5" '583.45'
$ws.Range("E5").Value = '  -0.48%  '

# --- Row 6: Solana ---
Set-TextValue "D6" '171.30'
$ws.Range("E6").Value = '  +2.01%  '

# --- Row 7: USDC ---
$ws.Range("E7").Value = '  -0.04%  '

# --- Row 8: XRP ---
$ws.Range("E8").Value = '  -0.93%  '

# --- Row 9: LidoStakedEther ---
Set-TextValue "D9" '2.485.35'
$ws.Range("E9").Value = '  -0.15%  '

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = '  +0.40%  '

# --- Row 12: Toncoin ---
$ws.Range("E12").Value = '  -0.47%  '

# --- Row 13: Cardano ---
Set-TextValue "D13" '0.332'
$ws.Range("E13").Value = '  -2.33%  '

# --- Row 15: (table row 15) ---
Set-TextValue "D15" '25.34'
$ws.Range("E15").Value = '  -2.56%  '

# --- Row 16 ---
Set-TextValue "D16" '67.060.87'
$ws.Range("E16").Value = '  +0.16%  '

# --- Row 17 ---
$ws.Range("E17").Value = '  -2.20%  '

# --- Row 18 ---
Set-TextValue "D18" '2.452.74'
$ws.Range("E18").Value = '  -0.65%  '

# --- Row 19 ---
$ws.Range("E19").Value = '  -6.14%  '

# --- Row 20: Uniswap ---
Set-TextValue "D20" '7.39'
$ws.Range("E20").Value = '  -5.17%  '

# --- Row 21 ---
Set-TextValue "D21" '347.95'
$ws.Range("E21").Value = '  -3.60%  '

# --- Row 22 ---
$ws.Range("E22").Value = '  -2.11%  '

# --- Row 23 ---
$ws.Range("E23").Value = '  -0.06%  '

# --- Row 24 ---
Set-TextValue "D24" '68.49'
$ws.Range("E24").Value = '  -3.18%  '

# --- Row 25 ---
$ws.Range("E25").Value = '  -4.81%  '

# --- Row 26 ---
$ws.Range("E26").Value = '  -3.63%  '

# --- Row 27 ---
$ws.Range("E27").Value = '  -2.02%  '

# --- Row 28: Binance-PegBSC-USD ---
Set-TextValue "D28" '1.00'
$ws.Range("E28").Value = '  +0.58%  '

# --- Row 30 ---
$ws.Range("D30").Value = '0.0₃0902'
$ws.Range("E30").Value = '  -3.86%  '

# --- Row 31 ---
Set-TextValue "D31" '509.18'
$ws.Range("E31").Value = '  +0.88%  '

# --- Row 32 ---
$ws.Range("E32").Value = '  -4.05%  '

# --- Row 33 ---
$ws.Range("E33").Value = '  -3.12%  '

# --- Row 34 ---
$ws.Range("E34").Value = '  -4.18%  '

# --- Row 35 ---
$ws.Range("E35").Value = '  +0.04%  '

# --- Row 36 ---
Set-TextValue "D36" '159.83'
$ws.Range("E36").Value = '  +0.35%  '

# --- Row 37 ---
Set-TextValue "D37" '0.116'
$ws.Range("E37").Value = '  -8.61%  '

# --- Row 38 ---
Set-TextValue "D38" '18.69'
$ws.Range("E38").Value = '  +0.68%  '

# --- Row 39 ---
$ws.Range("E39").Value = '  -5.16%  '

# --- Row 40 ---
Set-TextValue "D40" '1.33'
$ws.Range("E40").Value = '  -6.18%  '

# --- Row 41 ---
$ws.Range("E41").Value = '  -2.79%  '

# --- Row 42 ---
$ws.Range("E42").Value = '  -0.06%  '

# --- Row 43 ---
$ws.Range("E43").Value = '  -2.34%  '

# --- Row 44: RenderToken ---
$ws.Range("E44").Value = '  -3.28%  '

# --- Row 45: dogwifhat ---
$ws.Range("E45").Value = '  -5.10%  '

# --- Row 46: OKB ---
Set-TextValue "D46" '38.74'
$ws.Range("E46").Value = '  -1.52%  '

# --- Row 47: Aave ---
$ws.Range("E47").Value = '  +0.21%  '

# --- Row 48: ARBITRUM ---
$ws.Range("E48").Value = '  -4.76%  '

# --- Row 49: Filecoin ---
$ws.Range("E49").Value = '  -4.68%  '

# --- Rows 50 & 51: BabyDogeCoin/Cronos were swapped; fix so row 50 is
#     Cronos and row 51 is BabyDogeCoin (ranks in column A are unchanged).
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D50" '0.0732'
$ws.Range("E50").Value = '  -0.74%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0249'
$ws.Range("E51").Value = '  -6.39%  '
